# Auto-generated edit script: updates FFXIV leve profit figures (columns H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR worksheets, per the commit diff.

$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 272.0909
$ws.Range("I4").Value = 199.3
$ws.Range("K4").Value = 199.3
$ws.Range("M4").Value = -85.30000000000001
# Row 42
$ws.Range("H42").Value = 131.5
$ws.Range("J42").Value = 175.25
$ws.Range("L42").Value = 525.75
$ws.Range("N42").Value = -985.75
# Row 43
$ws.Range("H43").Value = 9806.532999999999
$ws.Range("J43").Value = 9899.923000000001
$ws.Range("L43").Value = 9899.923000000001
$ws.Range("N43").Value = -10037.923
# Row 107
$ws.Range("H107").Value = 333
$ws.Range("I107").Value = 229.22223
$ws.Range("K107").Value = 229.22223
$ws.Range("M107").Value = 1690.77777
# Row 125
$ws.Range("H125").Value = 2559.8333
$ws.Range("I125").Value = 2066.4
$ws.Range("J125").Value = 2912.2856
$ws.Range("K125").Value = 18597.6
$ws.Range("L125").Value = 26210.5704
$ws.Range("M125").Value = -16137.6
$ws.Range("N125").Value = -31130.5704
# Row 132
$ws.Range("H132").Value = 1059.2646
$ws.Range("I132").Value = 809.80646
$ws.Range("K132").Value = 2429.41938
$ws.Range("M132").Value = 100.5806199999997

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3662.2373
$ws.Range("I32").Value = 2677.2593
$ws.Range("K32").Value = 2677.2593
$ws.Range("M32").Value = -2390.2593
# Row 41
$ws.Range("H41").Value = 3053.4443
$ws.Range("I41").Value = 3053.4443
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 3053.4443
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -2639.4443
$ws.Range("N41").ClearContents()
# Row 61
$ws.Range("H61").Value = 4928.1787
$ws.Range("I61").Value = 3944.5
$ws.Range("K61").Value = 3944.5
$ws.Range("M61").Value = -3732.5
# Row 97
$ws.Range("H97").Value = 410.9091
$ws.Range("I97").Value = 372.45
$ws.Range("K97").Value = 372.45
$ws.Range("M97").Value = 123.55
# Row 136
$ws.Range("H136").Value = 4928.1787
$ws.Range("I136").Value = 3944.5
$ws.Range("K136").Value = 11833.5
$ws.Range("M136").Value = -9283.5
# Row 137
$ws.Range("H137").Value = 69998.8
$ws.Range("J137").Value = 69998.8
$ws.Range("L137").Value = 69998.8
$ws.Range("N137").Value = -80198.8
# Row 140
$ws.Range("H140").Value = 84209.5
$ws.Range("J140").Value = 84209.5
$ws.Range("L140").Value = 84209.5
$ws.Range("N140").Value = -94569.5

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 697.3333
$ws.Range("I80").Value = 716.3333
$ws.Range("J80").Value = 691
$ws.Range("K80").Value = 716.3333
$ws.Range("L80").Value = 691
$ws.Range("M80").Value = 281.6667
$ws.Range("N80").Value = -2687
# Row 83
$ws.Range("H83").Value = 697.3333
$ws.Range("I83").Value = 716.3333
$ws.Range("J83").Value = 691
$ws.Range("K83").Value = 3581.6665
$ws.Range("L83").Value = 3455
$ws.Range("M83").Value = 1410.3335
$ws.Range("N83").Value = -13439
# Row 100
$ws.Range("H100").Value = 15205.5
$ws.Range("J100").Value = 16266.6
$ws.Range("L100").Value = 16266.6
$ws.Range("N100").Value = -18430.6
# Row 129
$ws.Range("H129").Value = 68000
$ws.Range("J129").Value = 68000
$ws.Range("L129").Value = 68000
$ws.Range("N129").Value = -78000

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 5065.115
$ws.Range("I58").Value = 2442.0715
$ws.Range("J58").Value = 8125.3335
$ws.Range("K58").Value = 2442.0715
$ws.Range("L58").Value = 8125.3335
$ws.Range("M58").Value = -2239.0715
$ws.Range("N58").Value = -8531.333500000001
# Row 88
$ws.Range("H88").Value = 16515.4
$ws.Range("J88").Value = 17144.25
$ws.Range("L88").Value = 17144.25
$ws.Range("N88").Value = -17956.25
# Row 91
$ws.Range("H91").Value = 16515.4
$ws.Range("J91").Value = 17144.25
$ws.Range("L91").Value = 17144.25
$ws.Range("N91").Value = -19952.25
# Row 94
$ws.Range("H94").Value = 2084.8667
$ws.Range("J94").Value = 2620.7
$ws.Range("L94").Value = 2620.7
$ws.Range("N94").Value = -3522.7
# Row 95
$ws.Range("H95").Value = 20161.285
$ws.Range("J95").Value = 20161.285
$ws.Range("L95").Value = 20161.285
$ws.Range("N95").Value = -25653.285
# Row 132
$ws.Range("H132").Value = 3764.1667
$ws.Range("I132").Value = 2957.158
$ws.Range("K132").Value = 8871.474
$ws.Range("M132").Value = -6341.474
# Row 134
$ws.Range("H134").Value = 3292.7273
$ws.Range("I134").Value = 2690.4443
$ws.Range("K134").Value = 8071.3329
$ws.Range("M134").Value = -5536.3329
# Row 136
$ws.Range("H136").Value = 5065.115
$ws.Range("I136").Value = 2442.0715
$ws.Range("J136").Value = 8125.3335
$ws.Range("K136").Value = 7326.2145
$ws.Range("L136").Value = 24376.0005
$ws.Range("M136").Value = -4776.2145
$ws.Range("N136").Value = -29476.0005

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 173.24
$ws.Range("I2").Value = 168.92857
$ws.Range("J2").Value = 178.72728
$ws.Range("K2").Value = 1013.57142
$ws.Range("L2").Value = 1072.36368
$ws.Range("M2").Value = -900.57142
$ws.Range("N2").Value = -1298.36368
# Row 47
$ws.Range("H47").Value = 17749.834
$ws.Range("I47").Value = 20699.8
$ws.Range("J47").Value = 3000
$ws.Range("K47").Value = 62099.39999999999
$ws.Range("L47").Value = 9000
$ws.Range("M47").Value = -61668.39999999999
$ws.Range("N47").Value = -9862
# Row 113
$ws.Range("H113").Value = 52633268
$ws.Range("J113").Value = 83335140
$ws.Range("L113").Value = 250005420
$ws.Range("N113").Value = -250009760
# Row 114
$ws.Range("H114").Value = 4949
$ws.Range("J114").Value = 6282
$ws.Range("L114").Value = 18846
$ws.Range("N114").Value = -25354
# Row 134
$ws.Range("H134").Value = 22231624
$ws.Range("I134").Value = 10267.2
$ws.Range("J134").Value = 66674340
$ws.Range("K134").Value = 30801.6
$ws.Range("L134").Value = 200023020
$ws.Range("M134").Value = -25731.6
$ws.Range("N134").Value = -200033160
# Row 139
$ws.Range("H139").Value = 3694.7666
$ws.Range("J139").Value = 5694.3076
$ws.Range("L139").Value = 17082.9228
$ws.Range("N139").Value = -27362.9228
# Row 140
$ws.Range("H140").Value = 1275.56
$ws.Range("I140").Value = 1275.56
$ws.Range("K140").Value = 3826.68
$ws.Range("M140").Value = 1353.32

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 5315.4614
$ws.Range("I126").Value = 3174.75
$ws.Range("K126").Value = 9524.25
$ws.Range("M126").Value = -7054.25
# Row 133
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120
# Row 137
$ws.Range("H137").Value = 42447.5
# Row 140
$ws.Range("H140").Value = 69998.2
$ws.Range("J140").Value = 69998.2
$ws.Range("L140").Value = 69998.2
$ws.Range("N140").Value = -80358.2

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 10040
$ws.Range("I7").Value = 8401.556
$ws.Range("J7").Value = 12721.091
$ws.Range("K7").Value = 8401.556
$ws.Range("L7").Value = 12721.091
$ws.Range("M7").Value = -8289.556
$ws.Range("N7").Value = -12945.091
# Row 16
$ws.Range("H16").Value = 3980
$ws.Range("I16").Value = 2571.8
$ws.Range("J16").Value = 7500.5
$ws.Range("K16").Value = 2571.8
$ws.Range("L16").Value = 7500.5
$ws.Range("M16").Value = -2401.8
$ws.Range("N16").Value = -7840.5
# Row 18
$ws.Range("H18").Value = 500
$ws.Range("I18").Value = 500
$ws.Range("K18").Value = 500
$ws.Range("M18").Value = -328
# Row 40
$ws.Range("H40").Value = 16800
$ws.Range("I40").Value = 19711
$ws.Range("K40").Value = 19711
$ws.Range("M40").Value = -19575
# Row 46
$ws.Range("H46").Value = 2350.0444
$ws.Range("I46").Value = 1601.8096
$ws.Range("K46").Value = 1601.8096
$ws.Range("M46").Value = -1413.8096
# Row 55
$ws.Range("H55").Value = 1563908.6
$ws.Range("I55").Value = 2631954.8
$ws.Range("J55").Value = 2918.077
$ws.Range("K55").Value = 2631954.8
$ws.Range("L55").Value = 2918.077
$ws.Range("M55").Value = -2631781.8
$ws.Range("N55").Value = -3264.077
# Row 126
$ws.Range("H126").Value = 10040
$ws.Range("I126").Value = 8401.556
$ws.Range("J126").Value = 12721.091
$ws.Range("K126").Value = 25204.668
$ws.Range("L126").Value = 38163.273
$ws.Range("M126").Value = -22734.668
$ws.Range("N126").Value = -43103.273
# Row 132
$ws.Range("H132").Value = 16117.333
$ws.Range("I132").Value = 15913.286
$ws.Range("J132").Value = 16403
$ws.Range("K132").Value = 47739.858
$ws.Range("L132").Value = 49209
$ws.Range("M132").Value = -45209.858
$ws.Range("N132").Value = -54269
# Row 136
$ws.Range("H136").Value = 4323.52
$ws.Range("I136").Value = 3309.0571
$ws.Range("K136").Value = 9927.1713
$ws.Range("M136").Value = -7377.1713

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 19999.5
$ws.Range("I75").Value = 19999.5
$ws.Range("K75").Value = 19999.5
$ws.Range("M75").Value = -19063.5
# Row 78
$ws.Range("H78").Value = 19999.5
$ws.Range("I78").Value = 19999.5
$ws.Range("K78").Value = 59998.5
$ws.Range("M78").Value = -55318.5
# Row 113
$ws.Range("H113").Value = 403.875
$ws.Range("I113").Value = 298.58823
$ws.Range("J113").Value = 659.5714
$ws.Range("K113").Value = 895.76469
$ws.Range("L113").Value = 1978.7142
$ws.Range("M113").Value = 1274.23531
$ws.Range("N113").Value = -6318.7142
# Row 132
$ws.Range("H132").Value = 8742
$ws.Range("I132").Value = 4401.6665
$ws.Range("K132").Value = 13204.9995
$ws.Range("M132").Value = -10674.9995
